$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D2').Value = '36.724.88'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D3').Value = '1.963.54'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D5').Value = '244.92'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('D7').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D7').Value = '59.44'
$ws.Range('E7').Value = '  +2.94%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D9').Value = '0.374'
$ws.Range('E9').Value = '  +2.17%  '
$ws.Range('E10').Value = '  -2.03%  '
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('D12').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D12').Value = '22.47'
$ws.Range('E12').Value = '  +5.03%  '
$ws.Range('D13').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D13').Value = '2.252.96'
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D14').Value = '0.830'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('D15').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D15').Value = '13.76'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D17').Value = '1.974.86'
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('D18').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D18').Value = '36.579.90'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D19').Value = '69.98'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D20').Value = '0.0₃0861'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D21').Value = '5.09'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D22').Value = '229.16'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D24').Value = '2.45'
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('E25').Value = '  +3.28%  '
$ws.Range('E26').Value = '  +16.90%  '
$ws.Range('D27').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D27').Value = '9.27'
$ws.Range('E27').Value = '  +0.49%  '
$ws.Range('D28').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D28').Value = '160.83'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D29').Value = '19.42'
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('D31').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D32').Value = '4.74'
$ws.Range('E32').Value = '  +1.49%  '
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('D35').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D35').Value = '2.29'
$ws.Range('E35').Value = '  +7.20%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D37').Value = '3.41'
$ws.Range('E37').Value = '  +12.10%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D38').Value = '6.05'
$ws.Range('E38').Value = '  -2.37%  '
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D40').Value = '0.0986'
$ws.Range('E40').Value = '  +1.24%  '
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('E42').Value = '  +0.24%  '
$ws.Range('E43').Value = '  +1.89%  '
$ws.Range('D44').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D44').Value = '16.24'
$ws.Range('E44').Value = '  +0.95%  '
$ws.Range('D45').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D45').Value = '1.362.35'
$ws.Range('E45').Value = '  +0.83%  '
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('D47').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D47').Value = '88.01'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('D50').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D50').Value = '2.144.08'
$ws.Range('E50').Value = '  +1.34%  '
$ws.Range('D51').NumberFormat = '@'  # keep Price as text, matching source data
$ws.Range('D51').Value = '43.98'
$ws.Range('E51').Value = '  -2.67%  '
